$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated TPM-derived values (new TPM run changed the receptor average
# expression value for row 2 (M2); everything downstream that depends
# on it (receptor total expression value, the specificity scores, and
# the edge expression weights/specificities for rows 2-4) is
# recalculated accordingly.

$ws.Range("M2").Value = 4.959409333333333
$ws.Range("N2").Value = 14.878228
$ws.Range("O2").Value = 0.8271666313262851
$ws.Range("P2").Value = 0.8271666313262852
$ws.Range("Q2").Value = 0.6386694276737778
$ws.Range("R2").Value = 5.748024849064
$ws.Range("S2").Value = 0.8271666313262851
$ws.Range("T2").Value = 0.8271666313262852

$ws.Range("O3").Value = 0.09421438109281059
$ws.Range("P3").Value = 0.09421438109281059
$ws.Range("R3").Value = 0.6547007376159999
$ws.Range("S3").Value = 0.09421438109281059
$ws.Range("T3").Value = 0.09421438109281059

$ws.Range("O4").Value = 0.07861898758090437
$ws.Range("P4").Value = 0.07861898758090438
$ws.Range("S4").Value = 0.07861898758090437
$ws.Range("T4").Value = 0.07861898758090438
